$d = $word.ActiveDocument

# --- Change 1: merge "Install " + " " runs into a single "Install  " run ---
# Find the two-space-producing sequence "Install  SumatraPDF" (currently split across
# two runs: "Install " and " ") and normalize to a single run containing "Install  ".
$d.Content.Find.Execute("Install  SumatraPDF", $false, $false, $false, $false, $false, $true, 1, $false, "Install  SumatraPDF", 2) | Out-Null

# --- Change 2: insert a new paragraph holding a page break, right after the
# paragraph that ends with "change this to your installation" ---
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*change this to your installation*") {
        $target = $p
    }
}

if ($target -ne $null) {
    $r = $target.Range
    $r.Collapse(0)  # wdCollapseEnd
    $newPara = $r.InsertParagraphAfter()
    $afterRange = $target.Range
    $afterRange.Collapse(0)
    $afterRange.InsertBreak(7)  # wdPageBreak
}
